# ---------------------------------------------------------------------------
# Adds a new set of extracted-data rows (309-334) for the NTTHZ alloy family,
# sourced from 10.1016/j.actamat.2022.118602, to the "MiscSmallUploads" sheet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Reusable text blocks (kept as variables so the *first* time each unique
# string is written determines its position in the shared-strings table).
# ---------------------------------------------------------------------------
$txtComposition = "Hf15 Nb40 Ta25 Ti15 Zr5"
$txtNickname    = "NTTHZ"
$txtHL          = "cold rolled with 90% reduction and annealed at 900*C for 60min forming heterogeneous lamella (HL) structure"
$txtFG          = "cold rolled with 90% reduction and annealed at 1150*C for 60min forming fine grain (FG) structure"
$txtCG          = "cold rolled with 90% reduction and annealed at 1300*C for 30min forming coarse grain (CG) structure"
$txtMinDuct     = "minimum tensile ductility"
$txtAAMCR       = "AAM+CR"
$txtColdRolled  = "cold rolled with 90% reduction "
$txtFigS7a      = "S7a"
$txtFigF6a      = "F6a"
$txtFigS11a     = "S11a"
$txtDOI         = "10.1016/j.actamat.2022.118602"

$txtPhase       = "BCC"
$txtProcess     = "AAM+CR+A"
$txtYield       = "tensile yield stress"
$txtUTS         = "UTS"
$txtDuctility   = "tensile ductility"
$txtHardness    = "hardness"
$txtMethod      = "EXP"
$txtParams      = "strain rate 1e-3"
$txtPa          = "Pa"
$txtPct         = "%"
$txtAAM         = "AAM"

# ---------------------------------------------------------------------------
# Priming pass: touch each brand-new unique string for the first time in the
# exact order it first appears so the shared-strings table lines up with the
# target workbook (index 210 .. 221).
# ---------------------------------------------------------------------------
$ws.Range("B309").Value = $txtComposition   # -> shared string 210
$ws.Range("A309").Value = $txtNickname      # -> shared string 211
$ws.Range("E309").Value = $txtHL            # -> shared string 212
$ws.Range("E310").Value = $txtFG            # -> shared string 213
$ws.Range("E311").Value = $txtCG            # -> shared string 214
$ws.Range("F328").Value = $txtMinDuct       # -> shared string 215
$ws.Range("D331").Value = $txtAAMCR         # -> shared string 216
$ws.Range("E331").Value = $txtColdRolled    # -> shared string 217
$ws.Range("M330").Value = $txtFigS7a        # -> shared string 218
$ws.Range("M309").Value = $txtFigF6a        # -> shared string 219
$ws.Range("M332").Value = $txtFigS11a       # -> shared string 220
$ws.Range("N309").Value = $txtDOI           # -> shared string 221

# ---------------------------------------------------------------------------
# Row 309: NTTHZ, BCC, AAM+CR+A, HL structure, tensile yield stress @298K
# ---------------------------------------------------------------------------
$ws.Range("A309").Value = $txtNickname
$ws.Range("B309").Value = $txtComposition
$ws.Range("C309").Value = $txtPhase
$ws.Range("D309").Value = $txtProcess
$ws.Range("E309").Value = $txtHL
$ws.Range("F309").Value = $txtYield
$ws.Range("G309").Value = $txtMethod
$ws.Range("H309").Value = $txtParams
$ws.Range("I309").Value = 298
$ws.Range("J309").Value = 944444444.44444394
$ws.Range("L309").Value = $txtPa
$ws.Range("M309").Value = $txtFigF6a
$ws.Range("N309").Value = $txtDOI

# Row 310: FG, tensile yield stress @298K
$ws.Range("A310").Value = $txtNickname
$ws.Range("B310").Value = $txtComposition
$ws.Range("C310").Value = $txtPhase
$ws.Range("D310").Value = $txtProcess
$ws.Range("E310").Value = $txtFG
$ws.Range("F310").Value = $txtYield
$ws.Range("G310").Value = $txtMethod
$ws.Range("H310").Value = $txtParams
$ws.Range("I310").Value = 298
$ws.Range("J310").Value = 776353276.35327601
$ws.Range("L310").Value = $txtPa
$ws.Range("M310").Value = $txtFigF6a
$ws.Range("N310").Value = $txtDOI

# Row 311: CG, tensile yield stress @298K
$ws.Range("A311").Value = $txtNickname
$ws.Range("B311").Value = $txtComposition
$ws.Range("C311").Value = $txtPhase
$ws.Range("D311").Value = $txtProcess
$ws.Range("E311").Value = $txtCG
$ws.Range("F311").Value = $txtYield
$ws.Range("G311").Value = $txtMethod
$ws.Range("H311").Value = $txtParams
$ws.Range("I311").Value = 298
$ws.Range("J311").Value = 732193732.19373202
$ws.Range("L311").Value = $txtPa
$ws.Range("M311").Value = $txtFigF6a
$ws.Range("N311").Value = $txtDOI

# Row 312: HL, tensile yield stress @1073K
$ws.Range("A312").Value = $txtNickname
$ws.Range("B312").Value = $txtComposition
$ws.Range("C312").Value = $txtPhase
$ws.Range("D312").Value = $txtProcess
$ws.Range("E312").Value = $txtHL
$ws.Range("F312").Value = $txtYield
$ws.Range("G312").Value = $txtMethod
$ws.Range("H312").Value = $txtParams
$ws.Range("I312").Value = 1073
$ws.Range("J312").Value = 414529914.52991402
$ws.Range("L312").Value = $txtPa
$ws.Range("M312").Value = $txtFigF6a
$ws.Range("N312").Value = $txtDOI

# Row 313: HL, tensile yield stress @1173K
$ws.Range("A313").Value = $txtNickname
$ws.Range("B313").Value = $txtComposition
$ws.Range("C313").Value = $txtPhase
$ws.Range("D313").Value = $txtProcess
$ws.Range("E313").Value = $txtHL
$ws.Range("F313").Value = $txtYield
$ws.Range("G313").Value = $txtMethod
$ws.Range("H313").Value = $txtParams
$ws.Range("I313").Value = 1173
$ws.Range("J313").Value = 351851851.851852
$ws.Range("L313").Value = $txtPa
$ws.Range("M313").Value = $txtFigF6a
$ws.Range("N313").Value = $txtDOI

# Row 314: HL, tensile yield stress @1273K
$ws.Range("A314").Value = $txtNickname
$ws.Range("B314").Value = $txtComposition
$ws.Range("C314").Value = $txtPhase
$ws.Range("D314").Value = $txtProcess
$ws.Range("E314").Value = $txtHL
$ws.Range("F314").Value = $txtYield
$ws.Range("G314").Value = $txtMethod
$ws.Range("H314").Value = $txtParams
$ws.Range("I314").Value = 1273
$ws.Range("J314").Value = 243589743.58974299
$ws.Range("L314").Value = $txtPa
$ws.Range("M314").Value = $txtFigF6a
$ws.Range("N314").Value = $txtDOI

# Row 315: HL, tensile yield stress @1373K
$ws.Range("A315").Value = $txtNickname
$ws.Range("B315").Value = $txtComposition
$ws.Range("C315").Value = $txtPhase
$ws.Range("D315").Value = $txtProcess
$ws.Range("E315").Value = $txtHL
$ws.Range("F315").Value = $txtYield
$ws.Range("G315").Value = $txtMethod
$ws.Range("H315").Value = $txtParams
$ws.Range("I315").Value = 1373
$ws.Range("J315").Value = 142450142.450142
$ws.Range("L315").Value = $txtPa
$ws.Range("M315").Value = $txtFigF6a
$ws.Range("N315").Value = $txtDOI

# Row 316: HL, UTS @298K
$ws.Range("A316").Value = $txtNickname
$ws.Range("B316").Value = $txtComposition
$ws.Range("C316").Value = $txtPhase
$ws.Range("D316").Value = $txtProcess
$ws.Range("E316").Value = $txtHL
$ws.Range("F316").Value = $txtUTS
$ws.Range("G316").Value = $txtMethod
$ws.Range("H316").Value = $txtParams
$ws.Range("I316").Value = 298
$ws.Range("J316").Value = 1059829059.8290499
$ws.Range("L316").Value = $txtPa
$ws.Range("M316").Value = $txtFigF6a
$ws.Range("N316").Value = $txtDOI

# Row 317: FG, UTS @298K
$ws.Range("A317").Value = $txtNickname
$ws.Range("B317").Value = $txtComposition
$ws.Range("C317").Value = $txtPhase
$ws.Range("D317").Value = $txtProcess
$ws.Range("E317").Value = $txtFG
$ws.Range("F317").Value = $txtUTS
$ws.Range("G317").Value = $txtMethod
$ws.Range("H317").Value = $txtParams
$ws.Range("I317").Value = 298
$ws.Range("J317").Value = 913105413.10541296
$ws.Range("L317").Value = $txtPa
$ws.Range("M317").Value = $txtFigF6a
$ws.Range("N317").Value = $txtDOI

# Row 318: CG, UTS @298K
$ws.Range("A318").Value = $txtNickname
$ws.Range("B318").Value = $txtComposition
$ws.Range("C318").Value = $txtPhase
$ws.Range("D318").Value = $txtProcess
$ws.Range("E318").Value = $txtCG
$ws.Range("F318").Value = $txtUTS
$ws.Range("G318").Value = $txtMethod
$ws.Range("H318").Value = $txtParams
$ws.Range("I318").Value = 298
$ws.Range("J318").Value = 867521367.52136695
$ws.Range("L318").Value = $txtPa
$ws.Range("M318").Value = $txtFigF6a
$ws.Range("N318").Value = $txtDOI

# Row 319: HL, UTS @1073K
$ws.Range("A319").Value = $txtNickname
$ws.Range("B319").Value = $txtComposition
$ws.Range("C319").Value = $txtPhase
$ws.Range("D319").Value = $txtProcess
$ws.Range("E319").Value = $txtHL
$ws.Range("F319").Value = $txtUTS
$ws.Range("G319").Value = $txtMethod
$ws.Range("H319").Value = $txtParams
$ws.Range("I319").Value = 1073
$ws.Range("J319").Value = 538461538.46153796
$ws.Range("L319").Value = $txtPa
$ws.Range("M319").Value = $txtFigF6a
$ws.Range("N319").Value = $txtDOI

# Row 320: HL, UTS @1173K
$ws.Range("A320").Value = $txtNickname
$ws.Range("B320").Value = $txtComposition
$ws.Range("C320").Value = $txtPhase
$ws.Range("D320").Value = $txtProcess
$ws.Range("E320").Value = $txtHL
$ws.Range("F320").Value = $txtUTS
$ws.Range("G320").Value = $txtMethod
$ws.Range("H320").Value = $txtParams
$ws.Range("I320").Value = 1173
$ws.Range("J320").Value = 450142450.14244998
$ws.Range("L320").Value = $txtPa
$ws.Range("M320").Value = $txtFigF6a
$ws.Range("N320").Value = $txtDOI

# Row 321: HL, UTS @1273K
$ws.Range("A321").Value = $txtNickname
$ws.Range("B321").Value = $txtComposition
$ws.Range("C321").Value = $txtPhase
$ws.Range("D321").Value = $txtProcess
$ws.Range("E321").Value = $txtHL
$ws.Range("F321").Value = $txtUTS
$ws.Range("G321").Value = $txtMethod
$ws.Range("H321").Value = $txtParams
$ws.Range("I321").Value = 1273
$ws.Range("J321").Value = 313390313.39031303
$ws.Range("L321").Value = $txtPa
$ws.Range("M321").Value = $txtFigF6a
$ws.Range("N321").Value = $txtDOI

# Row 322: HL, UTS @1373K
$ws.Range("A322").Value = $txtNickname
$ws.Range("B322").Value = $txtComposition
$ws.Range("C322").Value = $txtPhase
$ws.Range("D322").Value = $txtProcess
$ws.Range("E322").Value = $txtHL
$ws.Range("F322").Value = $txtUTS
$ws.Range("G322").Value = $txtMethod
$ws.Range("H322").Value = $txtParams
$ws.Range("I322").Value = 1373
$ws.Range("J322").Value = 249287749.28774899
$ws.Range("L322").Value = $txtPa
$ws.Range("M322").Value = $txtFigF6a
$ws.Range("N322").Value = $txtDOI

# Row 323: HL, tensile ductility @298K
$ws.Range("A323").Value = $txtNickname
$ws.Range("B323").Value = $txtComposition
$ws.Range("C323").Value = $txtPhase
$ws.Range("D323").Value = $txtProcess
$ws.Range("E323").Value = $txtHL
$ws.Range("F323").Value = $txtDuctility
$ws.Range("G323").Value = $txtMethod
$ws.Range("H323").Value = $txtParams
$ws.Range("I323").Value = 298
$ws.Range("J323").Value = 13.988549618320601
$ws.Range("L323").Value = $txtPct
$ws.Range("M323").Value = $txtFigF6a
$ws.Range("N323").Value = $txtDOI

# Row 324: FG, tensile ductility @298K
$ws.Range("A324").Value = $txtNickname
$ws.Range("B324").Value = $txtComposition
$ws.Range("C324").Value = $txtPhase
$ws.Range("D324").Value = $txtProcess
$ws.Range("E324").Value = $txtFG
$ws.Range("F324").Value = $txtDuctility
$ws.Range("G324").Value = $txtMethod
$ws.Range("H324").Value = $txtParams
$ws.Range("I324").Value = 298
$ws.Range("J324").Value = 13.568702290076301
$ws.Range("L324").Value = $txtPct
$ws.Range("M324").Value = $txtFigF6a
$ws.Range("N324").Value = $txtDOI

# Row 325: CG, tensile ductility @298K
$ws.Range("A325").Value = $txtNickname
$ws.Range("B325").Value = $txtComposition
$ws.Range("C325").Value = $txtPhase
$ws.Range("D325").Value = $txtProcess
$ws.Range("E325").Value = $txtCG
$ws.Range("F325").Value = $txtDuctility
$ws.Range("G325").Value = $txtMethod
$ws.Range("H325").Value = $txtParams
$ws.Range("I325").Value = 298
$ws.Range("J325").Value = 15.706106870229
$ws.Range("L325").Value = $txtPct
$ws.Range("M325").Value = $txtFigF6a
$ws.Range("N325").Value = $txtDOI

# Row 326: HL, tensile ductility @1073K
$ws.Range("A326").Value = $txtNickname
$ws.Range("B326").Value = $txtComposition
$ws.Range("C326").Value = $txtPhase
$ws.Range("D326").Value = $txtProcess
$ws.Range("E326").Value = $txtHL
$ws.Range("F326").Value = $txtDuctility
$ws.Range("G326").Value = $txtMethod
$ws.Range("H326").Value = $txtParams
$ws.Range("I326").Value = 1073
$ws.Range("J326").Value = 9.2366412213740396
$ws.Range("L326").Value = $txtPct
$ws.Range("M326").Value = $txtFigF6a
$ws.Range("N326").Value = $txtDOI

# Row 327: HL, tensile ductility @1173K
$ws.Range("A327").Value = $txtNickname
$ws.Range("B327").Value = $txtComposition
$ws.Range("C327").Value = $txtPhase
$ws.Range("D327").Value = $txtProcess
$ws.Range("E327").Value = $txtHL
$ws.Range("F327").Value = $txtDuctility
$ws.Range("G327").Value = $txtMethod
$ws.Range("H327").Value = $txtParams
$ws.Range("I327").Value = 1173
$ws.Range("J327").Value = 10.7251908396946
$ws.Range("L327").Value = $txtPct
$ws.Range("M327").Value = $txtFigF6a
$ws.Range("N327").Value = $txtDOI

# Row 328: HL, minimum tensile ductility @1273K
$ws.Range("A328").Value = $txtNickname
$ws.Range("B328").Value = $txtComposition
$ws.Range("C328").Value = $txtPhase
$ws.Range("D328").Value = $txtProcess
$ws.Range("E328").Value = $txtHL
$ws.Range("F328").Value = $txtMinDuct
$ws.Range("G328").Value = $txtMethod
$ws.Range("H328").Value = $txtParams
$ws.Range("I328").Value = 1273
$ws.Range("J328").Value = 20
$ws.Range("L328").Value = $txtPct
$ws.Range("M328").Value = $txtFigF6a
$ws.Range("N328").Value = $txtDOI

# Row 329: HL, tensile ductility @1373K
$ws.Range("A329").Value = $txtNickname
$ws.Range("B329").Value = $txtComposition
$ws.Range("C329").Value = $txtPhase
$ws.Range("D329").Value = $txtProcess
$ws.Range("E329").Value = $txtHL
$ws.Range("F329").Value = $txtDuctility
$ws.Range("G329").Value = $txtMethod
$ws.Range("H329").Value = $txtParams
$ws.Range("I329").Value = 1373
$ws.Range("J329").Value = 9.1603053435114408
$ws.Range("L329").Value = $txtPct
$ws.Range("M329").Value = $txtFigF6a
$ws.Range("N329").Value = $txtDOI

# ---------------------------------------------------------------------------
# Row 330: as-cast hardness @298K, value computed from HV via formula
# ---------------------------------------------------------------------------
$ws.Range("A330").Value = $txtNickname
$ws.Range("B330").Value = $txtComposition
$ws.Range("C330").Value = $txtPhase
$ws.Range("D330").Value = $txtAAM
$ws.Range("F330").Value = $txtHardness
$ws.Range("G330").Value = $txtMethod
$ws.Range("I330").Value = 298
$ws.Range("P330").Value = 299
$ws.Range("L330").Value = $txtPa
$ws.Range("M330").Value = $txtFigS7a
$ws.Range("N330").Value = $txtDOI

# Row 331: cold-rolled hardness @298K, value computed from HV via formula
$ws.Range("A331").Value = $txtNickname
$ws.Range("B331").Value = $txtComposition
$ws.Range("C331").Value = $txtPhase
$ws.Range("D331").Value = $txtAAMCR
$ws.Range("E331").Value = $txtColdRolled
$ws.Range("F331").Value = $txtHardness
$ws.Range("G331").Value = $txtMethod
$ws.Range("I331").Value = 298
$ws.Range("P331").Value = 353
$ws.Range("L331").Value = $txtPa
$ws.Range("M331").Value = $txtFigS7a
$ws.Range("N331").Value = $txtDOI

# Shared formula J330:J331 = P*9807000 (Vickers hardness -> Pa)
$ws.Range("J330:J331").Formula = "=P330*9807000"

# ---------------------------------------------------------------------------
# Rows 332-334: cryogenic (77K) properties, HL condition
# ---------------------------------------------------------------------------
$ws.Range("A332").Value = $txtNickname
$ws.Range("B332").Value = $txtComposition
$ws.Range("C332").Value = $txtPhase
$ws.Range("D332").Value = $txtProcess
$ws.Range("E332").Value = $txtHL
$ws.Range("F332").Value = $txtYield
$ws.Range("G332").Value = $txtMethod
$ws.Range("H332").Value = $txtParams
$ws.Range("I332").Value = 77
$ws.Range("J332").Value = 1407000000
$ws.Range("J332").NumberFormat = "0.00E+00"
$ws.Range("L332").Value = $txtPa
$ws.Range("M332").Value = $txtFigS11a
$ws.Range("N332").Value = $txtDOI

$ws.Range("A333").Value = $txtNickname
$ws.Range("B333").Value = $txtComposition
$ws.Range("C333").Value = $txtPhase
$ws.Range("D333").Value = $txtProcess
$ws.Range("E333").Value = $txtHL
$ws.Range("F333").Value = $txtUTS
$ws.Range("G333").Value = $txtMethod
$ws.Range("H333").Value = $txtParams
$ws.Range("I333").Value = 77
$ws.Range("J333").Value = 1493000000
$ws.Range("J333").NumberFormat = "0.00E+00"
$ws.Range("L333").Value = $txtPa
$ws.Range("M333").Value = $txtFigS11a
$ws.Range("N333").Value = $txtDOI

$ws.Range("A334").Value = $txtNickname
$ws.Range("B334").Value = $txtComposition
$ws.Range("C334").Value = $txtPhase
$ws.Range("D334").Value = $txtProcess
$ws.Range("E334").Value = $txtHL
$ws.Range("F334").Value = $txtDuctility
$ws.Range("G334").Value = $txtMethod
$ws.Range("H334").Value = $txtParams
$ws.Range("I334").Value = 77
$ws.Range("J334").Value = 12.4
$ws.Range("L334").Value = $txtPct
$ws.Range("M334").Value = $txtFigS11a
$ws.Range("N334").Value = $txtDOI

# ---------------------------------------------------------------------------
# View-state: scroll the visible window down towards the newly-added rows
# and move the active selection, matching where the author's cursor ended up.
# ---------------------------------------------------------------------------
$ws.Range("A302").Select()
$excel.ActiveWindow.ScrollRow = 302
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N337").Select()
